# test_data.xlsx update
# 1. Set font in 'write_cell_content' func (red font for "未执行/blank" test-result cells)
# 2. Get data only from excel cell (account numbers are now read as text straight
#    from the sheet instead of being synthesized, and the login sheet's helper
#    formula in F2 is no longer needed)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("login")
$ws2 = $wb.Worksheets.Item("account")

# ---------------------------------------------------------------------------
# login sheet: drop the now-unused helper formula in F2
# ---------------------------------------------------------------------------
$ws1.Range("F2").ClearContents()

# ---------------------------------------------------------------------------
# account sheet: refresh the generated test account numbers/messages and
# store the account number as text (read verbatim from the cell rather than
# re-derived), plus mark the "测试结果" (test result) column
# ---------------------------------------------------------------------------

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
$ws2.Range("B2").Value = "jessicatest202007211"
Set-TextValue $ws2.Range("C2") "202007211"
$ws2.Range("E2").Value = 'Account "jessicatest202007211" was created.'
$ws2.Range("H2").ClearContents()
$ws2.Range("H2").Font.Color = 3158271

# Row 3
$ws2.Range("B3").Value = "jessicatest202007212"
Set-TextValue $ws2.Range("C3") "202007212"
$ws2.Range("E3").Value = 'Account "jessicatest202007212" was created.'
$ws2.Range("H3").Value = "NA"

# Row 4
$ws2.Range("B4").Value = "jessicatest202007213"
Set-TextValue $ws2.Range("C4") "202007213"
$ws2.Range("E4").Value = 'Account "jessicatest202007213" was created.'
$ws2.Range("H4").ClearContents()
$ws2.Range("H4").Font.Color = 3158271

# ---------------------------------------------------------------------------
# make "account" the active/selected sheet
# ---------------------------------------------------------------------------
$ws2.Activate()
